# Apply updates described by the commit "update infilled stock-recruit data time series"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "metadata": rename/expand definitions for C -> H, hr_pred -> H_cv,
# hr_pred_cv -> S_cv
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("metadata")

$meta.Range("A10").Value = "H"
$meta.Range("B10").Value = "Annual terminal harvest (i.e. catch)"

$meta.Range("A13").Value = "H_cv"
$meta.Range("B13").Value = "Coefficient of variation on harvest data. Historical (prior to 2011) Hucuktlis Sockeye harvest rate predictions were derived from a linear model. CV for these data is calculated as RMSE of the model residuals divided by the mean of the observed Hucuktlis Sockeye harvest rates that informed the model fit (i.e. the dependent variable). Harvest data for Somass and Hucuktlis post-2011 are assumed to be precise."

$meta.Range("A14").Value = "S_cv"
$meta.Range("B14").Value = "Coefficient of variation on spawner data. Currently based on ____"

# ---------------------------------------------------------------------------
# Sheet "S-R data": rename headers, then populate the new H_cv (L) and
# S_cv (M) columns across the whole data range.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("S-R data")

$ws.Range("I1").Value = "H"
$ws.Range("L1").Value = "H_cv"
$ws.Range("M1").Value = "S_cv"

# Rows 2-95: Somass/GCL/SPR records + the earliest HED rows with no
# existing CV data -> flat 0.05 / 0.05
for ($r = 2; $r -le 95; $r++) {
    $ws.Cells.Item($r, 12).Value = 0.05
    $ws.Cells.Item($r, 13).Value = 0.05
}

# Rows 96-129: HED rows that previously held the retrospective hr_pred /
# hr_pred_cv values -> replaced with the flat H_cv constant and a 0.2 S_cv
for ($r = 96; $r -le 129; $r++) {
    $ws.Cells.Item($r, 12).Value = 0.3899491183848233
    $ws.Cells.Item($r, 13).Value = 0.2
}

# Row 130 (2011): first year considered "precise" harvest data, but spawner
# data CV still elevated
$ws.Cells.Item(130, 12).Value = 0.05
$ws.Cells.Item(130, 13).Value = 0.2

# Rows 131-142 (2012-2023): precise harvest data, lower spawner CV
for ($r = 131; $r -le 142; $r++) {
    $ws.Cells.Item($r, 12).Value = 0.05
    $ws.Cells.Item($r, 13).Value = 0.1
}

# Rows 143-147 (1972-1976): earliest HED rows, use the flat H_cv constant
# and 0.2 S_cv
for ($r = 143; $r -le 147; $r++) {
    $ws.Cells.Item($r, 12).Value = 0.3899491183848233
    $ws.Cells.Item($r, 13).Value = 0.2
}
